# Autogenerated on Mon Feb 09 2015 03:30:35 GMT+0000 (Coordinated Universal Time)
# Adds the "Number of employees / Assets / Turnover" breakdown table
# (Micro/Small/Medium/Large) to the Guatemala Summary sheet, and moves the
# CIEN source citation further down (rows 29-30) to make room for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New table header (row 20) ---
$ws.Range("B20").Value = "Number of employees"
$ws.Range("C20").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D20").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B20:D20").Font.Bold = $true

# --- Micro row (row 21) - reuses the existing "Micro" text from row 11 ---
$ws.Range("A21").Value = "Micro"
$ws.Range("B21").Value = "0-10"

# --- Small row (row 22) ---
$ws.Range("A22").Value = "Small"
$ws.Range("B22").Value = "11-25"
$ws.Range("C22").Value = "<Q500,000 (Agexport)"

# --- Medium row (row 23) - replaces the old CIEN citation that lived here ---
$ws.Range("A23").Value = "Medium"
$ws.Range("A23").Font.Bold = $false
$ws.Range("B23").Value = "26-60"
$ws.Range("C23").Value = "<1,200,000 (Agexport)"

# --- Large row (row 24) - replaces the old long citation that lived here ---
$ws.Range("A24").Value = "Large"
$ws.Range("A24").Font.Italic = $false
$ws.Range("B24").Value = ">60"
$ws.Range("C24").Value = ">1,200,000 (Agexport)"

# --- Source citation, now moved down to rows 29-30 ---
$ws.Range("A29").Value = "CIEN"
$ws.Range("A29").Font.Bold = $true

$ws.Range("A30").Value = "CENTRO DE INVESTIGACIONES ECONÓMICAS NACIONALES, MICRO, PEQUENAS Y MEDIANAS EMPRESAS EN GUATEMALA. Available at http://www.mejoremosguate.org/cms/content/files/diagnosticos/economicos/Lineamientos_PYMES_05-05-2011.pdf"
$ws.Range("A30").Font.Italic = $true
